$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 172356
$ws.Range("C4").Value = 163148
$ws.Range("C7").Value = 5.34
$ws.Range("C8").Value = 65.92
